$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
Write-Host "ok"
